$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = 1.83
$ws.Range("K4").Value = 2.1
$ws.Range("L4").Value = 2.45
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 6.7
$ws.Range("O4").Value = 1.35
$ws.Range("T4").Value = 2.67
$ws.Range("U4").Value = 1.9
$ws.Range("W4").Value = 10.75
$ws.Range("AA4").Value = 45
$ws.Range("AC4").Value = 6.7
$ws.Range("AG4").Value = 6.2
$ws.Range("AH4").Value = 8
$ws.Range("AJ4").Value = 15
$ws.Range("AN4").Value = 5.9
$ws.Range("AT4").Value = 2.67
$ws.Range("AW4").Value = 3.65
$ws.Range("AX4").Value = 9.25
$ws.Range("AY4").Value = 19
